$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40 ---
$ws.Range("A40").Value = "I forget what day"
$ws.Range("B40").Value = "1700-I don’t remember"
$ws.Range("C40").Value = "Harry, Deon, Thuc"
$ws.Range("D40").Value = "Work on assignments"
$ws.Range("E40").Value = "Deon systematically figured out RuneLite’s architecture"
$ws.Range("G40").Value = "Today I have a monster diet. Android programming (from other courses) takes a lot of energy when doing for the first time."
$ws.Rows.Item(40).RowHeight = 49.25

# --- Row 41 ---
$ws.Range("A41").Value = "27 Feb 2020"
$ws.Range("B41").Value = "0100-0200"
$ws.Range("C41").Value = "Harry, Deon, Thuc"
$ws.Range("D41").Value = "Work on assignments"
$ws.Range("E41").Value = "We glossed over the document as the chaos from 261 kept us largely busy today"
$ws.Range("F41").Value = "Finding the architecture of a system is not easy, especially when people’s interpretations differ. Software is intangible and there is no realistic way (at least, of my knowledge) to verify whether nontrivial software can match their corresponding architecture."
$ws.Range("G41").Value = "Exhausted, once again, from the chaos of 261."
$ws.Rows.Item(41).RowHeight = 97

# --- Row 42 ---
$ws.Range("A42").Value = "27 Feb 2020"
$ws.Range("B42").Value = "1435-1525"
$ws.Range("C42").Value = "Harry"

$d42 = $ws.Range("D42")
$d42.Value = "Look at homework one last time."
$d42.Characters(1, 17).Font.Italic = $true
$d42.Characters(1, 17).Font.Color = 24832
$d42.Characters(18, 14).Font.Italic = $true
$d42.Characters(18, 14).Font.Strikethrough = $true
$d42.Characters(18, 14).Font.Color = 24832

$ws.Range("E42").Value = "Waiting for rest of team to show up. It’s 1600 and the others aren’t here."
$ws.Range("G42").Value = "Full. Just downed a whole Blaze Pizza (honestly, they’re not that big to begin with). My body is full of tomatoes now. Also, Sibelius’s Finlandia keeps my mood up."
$ws.Rows.Item(42).RowHeight = 61.15

# --- Row 43 ---
$ws.Range("A43").Value = "27 Feb 20020"
$ws.Range("B43").Value = "1620-1640"
$ws.Range("C43").Value = "Harry, Thuc"
$ws.Range("D43").Value = "Finalize the homework."
$ws.Range("E43").Value = "Pull request sent in."

# --- View state: scroll position & selection ---
$win = $excel.Windows.Item(1)
$win.ScrollRow = 40
$win.ScrollColumn = 3
$ws.Range("D43").Select()
